$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A, rows 2-24 currently hold numeric values 1..5 (one per question
# block of 5 rows). Convert them to text values with an "s" suffix
# (e.g. 1 -> "1s") so they are stored as shared strings instead of numbers.
for ($r = 2; $r -le 24; $r++) {
    $val = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = "$($val)s"
}

# Update the active selection to reflect where the user clicked after the
# edit (cell A25, just below the data).
$ws.Range("A25").Select()
